$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("area_pop_sum")

# Remove column C entirely (Density header, and C2/C3 values)
$ws.Range("C1:C3").Clear()

# Rename "Population" -> "population"
$ws.Range("A3").Value = "population"

# Add new row 4 with "density" label and the density value moved from C3
$ws.Range("A4").Value = "density"
$ws.Range("B4").Value = 1052.135608923543
